$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $cell = $ws.Range($rangeAddr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '27.247.41'
$ws.Range('D3').Value = '1.820.27'
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue 'D5' '313.39'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('E6').Value = '  +0.03%  '
Set-TextValue 'D7' '0.4461'
$ws.Range('E7').Value = '  -1.03%  '
Set-TextValue 'D8' '0.3768'
$ws.Range('E8').Value = '  +1.95%  '
Set-TextValue 'D9' '0.07392'
$ws.Range('E9').Value = '  +1.60%  '
Set-TextValue 'D10' '0.8784'
$ws.Range('E10').Value = '  +2.74%  '
Set-TextValue 'D11' '20.85'
$ws.Range('D12').Value = '1.815.42'
$ws.Range('E12').Value = '  -0.06%  '
Set-TextValue 'D13' '6.708'
$ws.Range('E13').Value = '  +0.94%  '
Set-TextValue 'D14' '5.413'
$ws.Range('E14').Value = '  +1.69%  '
Set-TextValue 'D15' '93.04'
$ws.Range('E15').Value = '  +0.82%  '
Set-TextValue 'D16' '0.07108'
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('E17').Value = '  +0.01%  '
Set-TextValue 'D18' '0.000008800'
$ws.Range('E18').Value = '  +0.16%  '
Set-TextValue 'D19' '1.001'
$ws.Range('E19').Value = '  +0.06%  '
Set-TextValue 'D20' '15.02'
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('D21').Value = '27.258.00'
$ws.Range('E21').Value = '  +1.10%  '
Set-TextValue 'D22' '5.350'
$ws.Range('E22').Value = '  +3.59%  '
$ws.Range('E23').Value = '  -0.09%  '
Set-TextValue 'D24' '1.960'
$ws.Range('E24').Value = '  -1.18%  '
Set-TextValue 'D25' '151.06'
$ws.Range('E25').Value = '  -0.34%  '
Set-TextValue 'D26' '2.296'
$ws.Range('E26').Value = '  +3.69%  '
Set-TextValue 'D27' '18.56'
$ws.Range('E27').Value = '  +0.86%  '
Set-TextValue 'D28' '5.343'
$ws.Range('E28').Value = '  +1.93%  '
Set-TextValue 'D29' '117.35'
$ws.Range('E29').Value = '  +0.84%  '
Set-TextValue 'D30' '0.08862'
$ws.Range('E30').Value = '  -0.04%  '
Set-TextValue 'D31' '0.7829'
$ws.Range('E31').Value = '  +4.25%  '
$ws.Range('E32').Value = '  +1.02%  '
Set-TextValue 'D33' '4.562'
$ws.Range('E33').Value = '  +2.78%  '
Set-TextValue 'D34' '2.906'
$ws.Range('E34').Value = '  -1.89%  '
$ws.Range('E35').Value = '  +0.04%  '
Set-TextValue 'D36' '1.109'
$ws.Range('E36').Value = '  +1.03%  '
Set-TextValue 'D37' '0.01965'
$ws.Range('E37').Value = '  -0.01%  '
Set-TextValue 'D38' '0.05257'
$ws.Range('E38').Value = '  +0.37%  '
Set-TextValue 'D39' '7.294'
$ws.Range('E39').Value = '  +1.66%  '
Set-TextValue 'D40' '0.5288'
$ws.Range('E40').Value = '  -0.35%  '
Set-TextValue 'D41' '2.858'
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range('E42').Value = '  -0.02%  '
Set-TextValue 'D43' '2.271'
$ws.Range('E43').Value = '  +14.90%  '
Set-TextValue 'D44' '8.599'
$ws.Range('E44').Value = '  +0.99%  '
Set-TextValue 'D45' '0.5025'
$ws.Range('E45').Value = '  -4.14%  '
Set-TextValue 'D46' '10.60'
$ws.Range('E46').Value = '  -0.19%  '
Set-TextValue 'D47' '104.89'
$ws.Range('E47').Value = '  -0.64%  '
Set-TextValue 'D48' '1.686'
$ws.Range('E48').Value = '  +1.19%  '
$ws.Range('E49').Value = '  +0.08%  '
Set-TextValue 'D50' '0.06384'
Set-TextValue 'D51' '66.04'
$ws.Range('E51').Value = '  +4.74%  '
